$d = $word.ActiveDocument

# Fix the split/misspelled "progrARTEs" (progr + ARTE + s, flagged by the
# spell-checker as a misspelling) so it reads correctly as "programas".
$d.Content.Find.Execute("progrARTEs", $true, $false, $false, $false, $false,
                         $true, 1, $false, "programas", 2)
